# Update the Pick & Place position data for U3 (Level Shifter) on row 18
# with the values from the newly generated (fixed) Gerber/placement files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "21.101mm"   # Mid X
$ws.Range("E18").Value = "21.274mm"   # Mid Y
$ws.Range("F18").Value = "14.732mm"   # Ref X
$ws.Range("G18").Value = "13.208mm"   # Ref Y
$ws.Range("H18").Value = "16.002mm"   # Pad X
$ws.Range("I18").Value = "14.924mm"   # Pad Y
